$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values (row 2 and row 3..5)
$ws.Range("B2").Value = 35227

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 19394

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3531

$ws.Range("B5").Value = 1327
